# Supprimer livraison.docx - apply commit "Modification des cas d'utilisations"
#
# Helper: Word's Find/Replace (and direct Range.Text assignment) that touches
# a run sometimes drops an adjacent *empty* trailing run in the same
# paragraph during re-serialization. When the diff requires that an empty
# run stay (or be introduced) right after the text run of a paragraph, call
# this helper right after editing that paragraph's text to (re)create it.
function Add-TrailingEmptyRun {
    param($para)
    $r = $para.Range
    $r.Collapse(0)              # wdCollapseEnd
    $r.InsertParagraphAfter()
    $endOfPara = $para.Range.End
    $markRange = $word.ActiveDocument.Range($endOfPara - 1, $endOfPara)
    $markRange.Delete()
}

$d = $word.ActiveDocument

# 1. Précondition: remove "(non) " before "valide"
$r = $d.Content
$r.Find.Execute(
    " Le système a calculé une tournée de livraison (non) valide et l’utilisateur a choisi de faire une modification sur le plan de livraison.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Le système a calculé une tournée de livraison valide et l’utilisateur a choisi de faire une modification sur le plan de livraison.",
    2) | Out-Null

# 2. "Le système demande ... choisir un (ou plusieur?) point ..." -> drop the
#    bold/red "(ou plusieur?)" remark entirely, merging the three runs into one.
$r = $d.Content
$r.Find.Execute(
    "Le système demande à l’utilisateur de choisir un (ou plusieur?) point de livraison à supprimer.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le système demande à l’utilisateur de choisir un point de livraison à supprimer.",
    2) | Out-Null

# 3. "L’utilisateur choisit un (ou plusieurs?) point ..." -> the bold/red run's
#    text becomes just a single space (instead of " (ou plusieurs?)"), and the
#    following plain run drops its own leading space accordingly.
$r = $d.Content
$r.Find.Execute(
    " (ou plusieurs?)", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$r.Text = " "

$r = $d.Content
$r.Find.Execute(
    " point de livraison à supprimer du plan de livraison et confirme son choix.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "point de livraison à supprimer du plan de livraison et confirme son choix.",
    2) | Out-Null

# 4. "Le système supprime le (les?) point ..." -> drop the bold/red "(les?)"
#    remark, merge the three runs into one, and extend the final sentence.
$r = $d.Content
$r.Find.Execute(
    "Le système supprime le (les?) point de livraison sélectionné du plan de livraison. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le système supprime le point de livraison sélectionné du plan de livraison et de la carte. ",
    2) | Out-Null

# 5. "Le système recalcule ..." -> new sentence about recomputing timings.
$r = $d.Content
$r.Find.Execute(
    "Le système recalcule une tournée de livraison et l’affiche dans la carte.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le système calcul le plus court chemin entre la livraison d’avance et la livraison d’après puis met à jour toutes les heures d’arrivées des livraisons suivantes.",
    2) | Out-Null

# 6. "L’utilisateur confirme le choix dans la nouvelle tournée calculée" ->
#    new wording, plus a new empty trailing run in the same paragraph.
$r = $d.Content
$r.Find.Execute(
    "L’utilisateur confirme le choix dans la nouvelle tournée calculée",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "L’utilisateur confirme le choix de la modification effectuée",
    2) | Out-Null
$paraConfirm = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("L’utilisateur confirme le choix de la modification effectuée")) {
        $paraConfirm = $d.Paragraphs($i)
        break
    }
}
Add-TrailingEmptyRun $paraConfirm

# 7. "2a. ... l’entrepôt parmis le (les) point de livraison à supprimer" ->
#    drop "(les) ".
$r = $d.Content
$r.Find.Execute(
    "2a. L’utilisateur décide de choisir l’entrepôt parmis le (les) point de livraison à supprimer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2a. L’utilisateur décide de choisir l’entrepôt parmis le point de livraison à supprimer",
    2) | Out-Null

# 8. "4a. La nouvelle tournée calculée ne rend pas une tournée valide" -> new
#    wording about time-slot constraints.
$r = $d.Content
$r.Find.Execute(
    "4a. La nouvelle tournée calculée ne rend pas une tournée valide",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "4a. La mise à jour des plages horaires provoque le non respect de la contrainte des plages horaires",
    2) | Out-Null

#    Next paragraph used to be just tab+tab then a bold/red "?" run; it
#    becomes tab+tab+plain sentence, merged into a single run, plus a new
#    empty trailing run.
$r = $d.Content
$r.Find.Execute(
    "^t^t?", $true, $false, $false, $false, $false, $true, 1, $false,
    "^t^tLe système met en surbrillance les plages horaires non valide",
    2) | Out-Null
$paraSurbrillance = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.IndexOf("Le système met en surbrillance les plages horaires non valide") -ge 0) {
        $paraSurbrillance = $d.Paragraphs($i)
        break
    }
}
Add-TrailingEmptyRun $paraSurbrillance

# 9. "Le système annule la suppression " -> append "d’une livraison"; keep the
#    paragraph's existing empty trailing run.
$r = $d.Content
$r.Find.Execute(
    "Le système annule la suppression ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le système annule la suppression d’une livraison",
    2) | Out-Null
$paraAnnule = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.IndexOf("Le système annule la suppression d’une livraison") -ge 0) {
        $paraAnnule = $d.Paragraphs($i)
        break
    }
}
Add-TrailingEmptyRun $paraAnnule
